# Apply crypto price/volume updates and row reordering per upstream data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "47.068.76"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +4.38%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.637.65"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +10.79%  "

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.02%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.19"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +6.88%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "104.55"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +10.71%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.611"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +9.81%  "

# Row 8
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.10%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.595"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +19.56%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.49"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +16.37%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "55.54"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +4.01%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0853"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +10.29%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "8.41"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +21.59%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.046.87"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +11.02%  "

# Row 15
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +3.25%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.657.98"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +11.90%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.943"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +15.10%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "15.36"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +10.34%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "47.684.14"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +5.82%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000103"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +11.51%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.22"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +7.01%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.81"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +12.33%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "72.48"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +9.57%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "270.70"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +13.76%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.10"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +13.18%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.20"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +17.72%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "30.31"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +44.88%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.23%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.08"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +1.48%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "10.72"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +13.01%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "40.09"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +6.69%  "

# Row 32
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +4.47%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.20"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +15.35%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.72"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -2.16%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.28"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +17.03%  "

# Row 36
$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0854"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +13.17%  "

# Row 37
$ws.Range("B37").Value = "WEMIXToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.88"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +6.26%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "151.93"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +3.28%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.124"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +11.20%  "

# Row 40
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +9.53%  "

# Row 41
$ws.Range("B41").Value = "EnergySwap"
$ws.Range("C41").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "23.59"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +57.54%  "

# Row 42
$ws.Range("B42").Value = "Celestia"
$ws.Range("C42").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "16.73"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +13.85%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.28"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +15.18%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.73"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +18.05%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0332"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +13.89%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.187.65"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +12.70%  "

# Row 47
$ws.Range("B47").Value = "FirstDigitalUSD"
$ws.Range("C47").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.999"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.02%  "

# Row 48
$ws.Range("B48").Value = "BitcoinSV"
$ws.Range("C48").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "95.53"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +7.21%  "

# Row 49
$ws.Range("B49").Value = "FraxShare"
$ws.Range("C49").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "10.15"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +20.23%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "114.49"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +15.62%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.84"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +7.11%  "
